$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '246.19'
Set-TextValue 'E2' '1.07%'
Set-TextValue 'D3' '29.57'
Set-TextValue 'E3' '-1.41%'
Set-TextValue 'D4' '5.155'
Set-TextValue 'E4' '0.03%'
Set-TextValue 'E5' '2.24%'
Set-TextValue 'D6' '6.652'
Set-TextValue 'E6' '1.75%'
Set-TextValue 'D7' '3.198'
Set-TextValue 'E7' '5.42%'
Set-TextValue 'D8' '0.8519'
Set-TextValue 'E8' '0.49%'
Set-TextValue 'D9' '0.8653'
Set-TextValue 'E9' '0.51%'
Set-TextValue 'B10' 'One'
Set-TextValue 'C10' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D10' '0.01030'
Set-TextValue 'E10' '1,613.88%'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1377'
Set-TextValue 'E11' '2.07%'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.07091'
Set-TextValue 'E12' '2.59%'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03265'
Set-TextValue 'E13' '13.01%'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09374'
Set-TextValue 'E14' '-0.06%'
Set-TextValue 'B15' 'BitForexToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001537'
Set-TextValue 'E15' '0.72%'
Set-TextValue 'D16' '0.006145'
Set-TextValue 'E16' '1.14%'
Set-TextValue 'E17' '-0.61%'
Set-TextValue 'E18' '-0.94%'
Set-TextValue 'D20' '0.03366'
Set-TextValue 'E20' '0.42%'
Set-TextValue 'D21' '0.1282'
Set-TextValue 'E21' '-1.64%'
Set-TextValue 'D22' '3.322'
Set-TextValue 'E22' '-8.57%'
Set-TextValue 'E23' '-0.76%'
Set-TextValue 'D24' '0.1380'
Set-TextValue 'E24' '0.50%'
Set-TextValue 'E25' '1.20%'
Set-TextValue 'D26' '0.004143'
Set-TextValue 'E26' '-6.64%'
Set-TextValue 'D27' '0.0001209'
Set-TextValue 'E27' '2.52%'
Set-TextValue 'D28' '0.0001444'
Set-TextValue 'E28' '3.85%'
Set-TextValue 'D40' '0.03747'
Set-TextValue 'E40' '0.02%'
Set-TextValue 'D41' '0.005816'
Set-TextValue 'E41' '-0.54%'
Set-TextValue 'D42' '0.1070'
Set-TextValue 'E42' '1.39%'
Set-TextValue 'D43' '0.002199'
Set-TextValue 'E43' '-4.37%'
Set-TextValue 'D44' '0.009182'
Set-TextValue 'E44' '-1.12%'
Set-TextValue 'D45' '0.00005300'
Set-TextValue 'E45' '3.78%'
Set-TextValue 'D46' '0.00000000750'
Set-TextValue 'E46' '-0.02%'
Set-TextValue 'D47' '0.05796'
Set-TextValue 'E47' '-42.01%'
Set-TextValue 'D48' '0.002175'
Set-TextValue 'E48' '-21.58%'
Set-TextValue 'D49' '0.00002099'
Set-TextValue 'E49' '-0.02%'
Set-TextValue 'D50' '0.0001999'
Set-TextValue 'E50' '-0.02%'
